# "added names we forgot" - fill in the blank rows in the people roster
# and fix the Name column header value for Sabien Jarmin to include her
# handle.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("people")

# New team members that were missing from the roster
$ws.Range("D16").Value = 'Shane "Vhespir"'
$ws.Range("E16").Value = "Composer"

$ws.Range("D17").Value = 'Logan "Ryver" Fairbairn'
$ws.Range("E17").Value = "Character artist"

# Sabien's name gained her handle "ambid"
$ws.Range("D2").Value = 'Sabien "ambid" Jarmin'

$ws.Range("D18").Value = '"esvento"'
$ws.Range("E18").Value = "Character artist"

$ws.Range("D19").Value = "catarina"
$ws.Range("E19").Value = "Concept artist"

$ws.Range("D20").Value = 'Sophie "Spoot"'
$ws.Range("E20").Value = "Concept artist"

# Update viewport / selection to match the author's final cursor position
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("E19").Select()
